# Daily attendance processing - 2026-01-07 20:38:34
# Swap the order of "Recorded By" entries from "dnasr281@gmail.com, System"
# to "System, dnasr281@gmail.com" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -eq $oldText) {
        $cell.Value = $newText
        $changed = $changed + 1
    }
}

Write-Host "Updated cells:" $changed
